$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Roll the quarterly window forward by one quarter ---
# Drop the oldest quarter (1399/06), shift D:L values left into C:K->D:L,
# and populate the newest quarter (1401/12) into column M.

# Row 8: period-ending labels
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# Row 9: statement publish dates
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-04-18 (8)"
$ws.Range("F9").Value = "1401-04-28 (2)"
$ws.Range("G9").Value = "1401-08-29 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-23 (7)"
$ws.Range("K9").Value = "1401-08-29 (2)"

# J9, L9 and M9 hold bare "yyyy-mm-dd" text (no "(n)" suffix), which Excel's
# smart-typing would otherwise silently coerce into a date serial number.
# Force them to Text first, then re-copy the known-good text formatting
# (border/fill/font/alignment) from a sibling cell in the same row so the
# cell's style index stays identical to its neighbours, same as the source file.
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "1401-04-28"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "1402-02-23"

$ws.Range("K9").Copy() | Out-Null
$ws.Range("J9").PasteSpecial(-4122) | Out-Null
$ws.Range("K9").Copy() | Out-Null
$ws.Range("L9").PasteSpecial(-4122) | Out-Null
$ws.Range("K9").Copy() | Out-Null
$ws.Range("M9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows 11-27: financial statement figures (read_price algorithm refresh)
# Row 11: Sales
$ws.Range("D11").Value = 198132171
$ws.Range("E11").Value = 318788584
$ws.Range("F11").Value = 246141824
$ws.Range("G11").Value = 262090687
$ws.Range("H11").Value = 425188986
$ws.Range("I11").Value = 478985412
$ws.Range("J11").Value = 580668052
$ws.Range("K11").Value = 586359486
$ws.Range("L11").Value = 560360321
$ws.Range("M11").Value = 473728011

# Row 12: Cost of goods sold
$ws.Range("D12").Value = -177605018
$ws.Range("E12").Value = -258434779
$ws.Range("F12").Value = -221605924
$ws.Range("G12").Value = -232828201
$ws.Range("H12").Value = -382791344
$ws.Range("I12").Value = -439010858
$ws.Range("J12").Value = -464864867
$ws.Range("K12").Value = -532138948
$ws.Range("L12").Value = -512261439
$ws.Range("M12").Value = -430005346

# Row 13: Gross profit (loss)
$ws.Range("D13").Value = 20527153
$ws.Range("E13").Value = 60353805
$ws.Range("F13").Value = 24535900
$ws.Range("G13").Value = 29262486
$ws.Range("H13").Value = 42397642
$ws.Range("I13").Value = 39974554
$ws.Range("J13").Value = 115803185
$ws.Range("K13").Value = 54220538
$ws.Range("L13").Value = 48098882
$ws.Range("M13").Value = 43722665

# Row 14: General, administrative & organizational expenses
$ws.Range("D14").Value = -373301
$ws.Range("E14").Value = -8308989
$ws.Range("F14").Value = -1041064
$ws.Range("G14").Value = -1094573
$ws.Range("H14").Value = -2142329
$ws.Range("I14").Value = -5989652
$ws.Range("J14").Value = -2031322
$ws.Range("K14").Value = -3136576
$ws.Range("L14").Value = -606390
$ws.Range("M14").Value = -8239437

# Row 15: Impairment of receivables (exceptional expense)
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16: Other operating income (expense), net
$ws.Range("D16").Value = 464478
$ws.Range("E16").Value = 10556155
$ws.Range("F16").Value = 510922
$ws.Range("G16").Value = 905912
$ws.Range("H16").Value = 9092114
$ws.Range("I16").Value = 4282304
$ws.Range("J16").Value = 1766253
$ws.Range("K16").Value = 1138873
$ws.Range("L16").Value = 1400310
$ws.Range("M16").Value = 1712834

# Row 17: Operating profit (loss)
$ws.Range("D17").Value = 20618330
$ws.Range("E17").Value = 62600971
$ws.Range("F17").Value = 24005758
$ws.Range("G17").Value = 29073825
$ws.Range("H17").Value = 49347427
$ws.Range("I17").Value = 38267206
$ws.Range("J17").Value = 115538116
$ws.Range("K17").Value = 52222835
$ws.Range("L17").Value = 48892802
$ws.Range("M17").Value = 37196062

# Row 18: Finance costs
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 35004
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = -499315
$ws.Range("L18").Value = -1337205
$ws.Range("M18").Value = -2008024

# Row 19: Other non-operating income and expense, net
$ws.Range("D19").Value = 4513486
$ws.Range("E19").Value = 5740537
$ws.Range("F19").Value = 6429396
$ws.Range("G19").Value = 9729790
$ws.Range("H19").Value = 8585424
$ws.Range("I19").Value = 10260110
$ws.Range("J19").Value = 11241737
$ws.Range("K19").Value = 21534224
$ws.Range("L19").Value = 17809916
$ws.Range("M19").Value = 21949859

# Row 20: Net profit (loss) from continuing operations before tax
$ws.Range("D20").Value = 25131816
$ws.Range("E20").Value = 68376512
$ws.Range("F20").Value = 30435154
$ws.Range("G20").Value = 38803615
$ws.Range("H20").Value = 57932851
$ws.Range("I20").Value = 48527316
$ws.Range("J20").Value = 126779853
$ws.Range("K20").Value = 73257744
$ws.Range("L20").Value = 65365513
$ws.Range("M20").Value = 57137897

# Row 21: Tax
$ws.Range("D21").Value = -4141147
$ws.Range("E21").Value = -3369672
$ws.Range("F21").Value = -4726212
$ws.Range("G21").Value = -7772188
$ws.Range("H21").Value = -9835519
$ws.Range("I21").Value = 534826
$ws.Range("J21").Value = -25008565
$ws.Range("K21").Value = -11860790
$ws.Range("L21").Value = 506389
$ws.Range("M21").Value = 13853751

# Row 22: Net profit (loss) from continuing operations
$ws.Range("D22").Value = 20990669
$ws.Range("E22").Value = 65006840
$ws.Range("F22").Value = 25708942
$ws.Range("G22").Value = 31031427
$ws.Range("H22").Value = 48097332
$ws.Range("I22").Value = 49062142
$ws.Range("J22").Value = 101771288
$ws.Range("K22").Value = 61396954
$ws.Range("L22").Value = 65871902
$ws.Range("M22").Value = 70991648

# Row 23: Profit (loss) from discontinued operations, net of tax
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24: Net profit (loss)
$ws.Range("D24").Value = 20990669
$ws.Range("E24").Value = 65006840
$ws.Range("F24").Value = 25708942
$ws.Range("G24").Value = 31031427
$ws.Range("H24").Value = 48097332
$ws.Range("I24").Value = 49062142
$ws.Range("J24").Value = 101771288
$ws.Range("K24").Value = 61396954
$ws.Range("L24").Value = 65871902
$ws.Range("M24").Value = 70991648

# Row 25: Earnings per share after tax
$ws.Range("D25").Value = 123
$ws.Range("E25").Value = 382
$ws.Range("F25").Value = 151
$ws.Range("G25").Value = 183
$ws.Range("H25").Value = 175
$ws.Range("I25").Value = 178
$ws.Range("J25").Value = 599
$ws.Range("K25").Value = 361
$ws.Range("L25").Value = 240
$ws.Range("M25").Value = 258

# Row 26: Capital
$ws.Range("D26").Value = 170000000
$ws.Range("E26").Value = 170000000
$ws.Range("F26").Value = 170000000
$ws.Range("G26").Value = 170000000
$ws.Range("H26").Value = 275000000
$ws.Range("I26").Value = 275000000
$ws.Range("J26").Value = 170000000
$ws.Range("K26").Value = 170000000
$ws.Range("L26").Value = 275000000
$ws.Range("M26").Value = 275000000

# Row 27: Earnings per share based on latest capital
$ws.Range("D27").Value = 76
$ws.Range("E27").Value = 236
$ws.Range("F27").Value = 93
$ws.Range("G27").Value = 113
$ws.Range("H27").Value = 175
$ws.Range("I27").Value = 178
$ws.Range("J27").Value = 370
$ws.Range("K27").Value = 223
$ws.Range("L27").Value = 240
$ws.Range("M27").Value = 258

# --- Column widths: the "year-end quarter" (Q4) columns get the wider 31-char width ---
# Before: F (1399/12) and J (1400/12) were the wide columns.
# After: E (1399/12), I (1400/12) and the new M (1401/12) are the wide columns.
$wNormal = $ws.Columns.Item(4).ColumnWidth   # current width of column D (the "29" class)
$wWide   = $ws.Columns.Item(6).ColumnWidth   # current width of column F (the "31" class)

$ws.Columns.Item(4).ColumnWidth  = $wNormal  # D
$ws.Columns.Item(5).ColumnWidth  = $wWide    # E
$ws.Columns.Item(6).ColumnWidth  = $wNormal  # F
$ws.Columns.Item(7).ColumnWidth  = $wNormal  # G
$ws.Columns.Item(8).ColumnWidth  = $wNormal  # H
$ws.Columns.Item(9).ColumnWidth  = $wWide    # I
$ws.Columns.Item(10).ColumnWidth = $wNormal  # J
$ws.Columns.Item(11).ColumnWidth = $wNormal  # K
$ws.Columns.Item(12).ColumnWidth = $wNormal  # L
$ws.Columns.Item(13).ColumnWidth = $wWide    # M
